# Fix bugs in holdAndWin data table: rows 2-25 (columns A-F) were scrambled.
# This restores the correct symbol-row ordering and reel-weight values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2 = @(201, 9, 30, 15, 45, 30)
    3 = @(1201, 2, 10, 10, 10, 10)
    4 = @(1202, 2, 10, 10, 10, 10)
    5 = @(101, 9, 30, 15, 60, 15)
    6 = @(401, 9, 48, 67, 75, 45)
    7 = @(701, 3, 90, 45, 97, 15)
    8 = @(801, 3, 67, 65, 52, 45)
    9 = @(1203, 3, 15, 15, 15, 15)
    10 = @(901, 16, 15, 45, 60, 60)
    11 = @(501, 9, 52, 30, 75, 45)
    12 = @(301, 6, 45, 30, 60, 45)
    13 = @(601, 9, 60, 67, 60, 42)
    14 = @(902, 1, 0, 0, 0, 0)
    15 = @(1001, 18, 30, 75, 60, 72)
    16 = @(802, 0, 4, 5, 4, 0)
    17 = @(2, 0, 2, 2, 2, 2)
    18 = @(1101, 0, 15, 30, 30, 0)
    19 = @(3, 0, 3, 3, 3, 3)
    20 = @(502, 0, 4, 0, 0, 0)
    21 = @(1, 0, 2, 2, 2, 2)
    22 = @(402, 0, 0, 4, 0, 0)
    23 = @(602, 0, 0, 4, 0, 9)
    24 = @(702, 0, 0, 0, 4, 0)
    25 = @(1002, 0, 0, 0, 0, 9)
}

foreach ($row in $values.Keys) {
    $rowVals = $values[$row]
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($row, $col).Value = $rowVals[$col - 1]
    }
}
